# Updated cryptos list on Fri Nov  3 09:01:10 UTC 2023 with GitHub Actions
#
# Applies the refreshed price / volume(1h) figures to the crypto table on
# Sheet1, including the Polygon / WrappedEther rank swap between rows 14
# and 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# The "Price" column (D) holds numeric-looking strings (e.g. "230.78",
# "0.0681") that must stay plain text, exactly like the source data - not
# get auto-coerced into Excel numbers. Force text storage via NumberFormat
# "@" before writing, then restore the "Normal" style afterwards so no
# stray cell-format diff is left behind.
function Set-PriceCell($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceCell "D2" "34.721.86"
Set-Cell "E2" "  -2.05%  "

# Row 3 - Ethereum
Set-PriceCell "D3" "1.809.21"
Set-Cell "E3" "  -1.60%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.38%  "

# Row 5 - BNB
Set-PriceCell "D5" "230.78"
Set-Cell "E5" "  -0.18%  "

# Row 6 - XRP
Set-PriceCell "D6" "0.608"
Set-Cell "E6" "  -0.59%  "

# Row 7 - USDC
Set-Cell "E7" "  +0.48%  "

# Row 8 - Solana
Set-PriceCell "D8" "39.44"
Set-Cell "E8" "  -10.44%  "

# Row 9 - Cardano
Set-PriceCell "D9" "0.325"
Set-Cell "E9" "  +4.95%  "

# Row 10 - Dogecoin
Set-PriceCell "D10" "0.0681"
Set-Cell "E10" "  -3.64%  "

# Row 11 - TRON
Set-PriceCell "D11" "0.0993"
Set-Cell "E11" "  -1.55%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceCell "D12" "2.071.83"
Set-Cell "E12" "  -1.51%  "

# Row 13 - Chainlink
Set-PriceCell "D13" "11.24"
Set-Cell "E13" "  -0.44%  "

# Rows 14/15 - Polygon and WrappedEther swap rank positions
Set-Cell "B14" "Polygon"
Set-Cell "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-PriceCell "D14" "0.665"
Set-Cell "E14" "  -1.40%  "

Set-Cell "B15" "WrappedEther"
Set-Cell "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PriceCell "D15" "1.809.48"
Set-Cell "E15" "  -1.49%  "

# Row 16 - Polkadot
Set-PriceCell "D16" "4.59"
Set-Cell "E16" "  -2.29%  "

# Row 17 - WrappedBTC
Set-PriceCell "D17" "34.698.55"
Set-Cell "E17" "  -2.06%  "

# Row 18 - Litecoin
Set-PriceCell "D18" "69.48"
Set-Cell "E18" "  -0.77%  "

# Row 19 - ShibaInu
Set-PriceCell "D19" "0.0₃0785"
Set-Cell "E19" "  -2.13%  "

# Row 20 - BitcoinCash
Set-PriceCell "D20" "240.39"
Set-Cell "E20" "  -1.65%  "

# Row 21 - Avalanche
Set-PriceCell "D21" "11.92"
Set-Cell "E21" "  -1.61%  "

# Row 22 - Uniswap
Set-PriceCell "D22" "4.68"
Set-Cell "E22" "  -0.65%  "

# Row 23 - Dai
Set-Cell "E23" "  +0.55%  "

# Row 24 - Toncoin
Set-Cell "E24" "  +1.67%  "

# Row 25 - Monero
Set-PriceCell "D25" "172.02"
Set-Cell "E25" "  +1.26%  "

# Row 26 - Cosmos
Set-PriceCell "D26" "7.78"
Set-Cell "E26" "  -1.77%  "

# Row 27 - EthereumClassic
Set-PriceCell "D27" "17.23"
Set-Cell "E27" "  -3.03%  "

# Row 28 - Stellar
Set-PriceCell "D28" "0.121"
Set-Cell "E28" "  -0.12%  "

# Row 29 - PancakeSwap
Set-Cell "E29" "  +0.01%  "

# Row 30 - BinanceUSD
Set-Cell "E30" "  +0.41%  "

# Row 31 - Filecoin
Set-Cell "E31" "  +3.40%  "

# Row 32 - Hedera
Set-PriceCell "D32" "0.0547"
Set-Cell "E32" "  -1.22%  "

# Row 33 - InternetComputer(DFINITY)
Set-PriceCell "D33" "3.94"
Set-Cell "E33" "  -3.80%  "

# Row 34 - TrustWalletToken
Set-PriceCell "D34" "1.28"
Set-Cell "E34" "  +16.12%  "

# Row 35 - LidoDAOToken
Set-PriceCell "D35" "1.79"
Set-Cell "E35" "  -3.77%  "

# Row 36 - ImmutableX
Set-PriceCell "D36" "0.699"
Set-Cell "E36" "  +1.81%  "

# Row 37 - Aave
Set-PriceCell "D37" "91.78"
Set-Cell "E37" "  -4.99%  "

# Row 38 - WEMIXToken
Set-Cell "E38" "  +4.39%  "

# Row 39 - Maker
Set-PriceCell "D39" "1.326.90"
Set-Cell "E39" "  -1.65%  "

# Row 40 - VeChain
Set-PriceCell "D40" "0.0192"
Set-Cell "E40" "  -1.51%  "

# Row 41 - HuobiToken
Set-Cell "E41" "  +0.55%  "

# Row 42 - ARBITRUM
Set-PriceCell "D42" "0.964"
Set-Cell "E42" "  -5.09%  "

# Row 43 - InjectiveProtocol
Set-PriceCell "D43" "14.34"
Set-Cell "E43" "  -8.74%  "

# Row 44 - RenderToken
Set-Cell "E44" "  -9.27%  "

# Row 45 - MXToken
Set-PriceCell "D45" "2.70"
Set-Cell "E45" "  -3.90%  "

# Row 46 - FraxShare
Set-PriceCell "D46" "6.24"
Set-Cell "E46" "  -0.37%  "

# Row 47 - Kaspa
Set-PriceCell "D47" "0.0513"
Set-Cell "E47" "  -1.23%  "

# Row 48 - RocketPoolETH
Set-PriceCell "D48" "1.996.39"
Set-Cell "E48" "  -0.42%  "

# Row 49 - PaxDollar
Set-Cell "E49" "  +0.47%  "

# Row 50 - Cronos
Set-PriceCell "D50" "0.0666"
Set-Cell "E50" "  +6.03%  "

# Row 51 - Quant
Set-PriceCell "D51" "98.13"
Set-Cell "E51" "  -5.12%  "
